$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recipes")

# --- Row 1: swap header labels in A1/B1 ---
$ws.Range("A1").Value = "ingredients"
$ws.Range("B1").Value = "steps"

# --- Row 2: swap the steps/ingredients text between A2 and B2 ---
$stepsText = "Measure 1 cup of sugar into a mixing bowl.`nMeasure 1 cup of butter into the same mixing bowl.`nMix the sugar and butter together until combined.`nAdd eggs to the same bowl and mix until combined.`nMeasure 2 teaspoons of vanilla extract into the same mixing bowl.`nMix the vanilla extract until combined. `nMeasure 1.5 cups of all purpose flour into a new mixing bowl.`nMeasure 1.75 teaspoons of baking powder into bowl with the flour.`nBake.`n"
$ingredientsText = "1 cup white sugar, 0.5 cups unsalted butter, 2 teaspoons vanilla extract, 1.5 cups all purpose flour, 1.75 teaspoons baking powder, 0.5 cups milk"

$ws.Range("A2").Value = $ingredientsText
$ws.Range("B2").Value = $stepsText

# --- Row 2: update classification value ---
$ws.Range("E2").Value = "1, 1"

# --- Row 3: add new recipe entry ---
$ws.Range("A3").Value = "Test"
$ws.Range("B3").Value = "Test`n"
$ws.Range("C3").Value = "Test"
$ws.Range("D3").Value = 21

# E3 must stay text "0" (not numeric 0) - force Text format before
# assigning, then drop back to the default "Normal" style so we don't
# leave a stray numeric-style cell behind.
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0"
$ws.Range("E3").Style = "Normal"

$wb.Save()
